$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (I1, J1) with same style as H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows for columns I (I0) and J (IF)
$data = @{
    2 = @(8, 9)
    3 = @(9, 9)
    4 = @(9, 9)
    5 = @(8, 8)
    6 = @(7, 7)
    7 = @(6, 7)
    8 = @(6, 7)
    9 = @(9, 9)
    10 = @(8, 8)
    11 = @(3, 4)
    12 = @(7, 8)
    13 = @(6, 6)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(7, 7)
    17 = @(5, 6)
    18 = @(6, 6)
    19 = @(6, 6)
    20 = @(6, 7)
    21 = @(7, 7)
    22 = @(6, 7)
    23 = @(8, 8)
    24 = @(8, 8)
    25 = @(7, 7)
    26 = @(5, 7)
    27 = @(8, 9)
    28 = @(8, 9)
    29 = @(8, 8)
    30 = @(5, 5)
    31 = @(7, 7)
    32 = @(8, 9)
    33 = @(6, 6)
    34 = @(8, 8)
    35 = @(7, 8)
    36 = @(10, 10)
    37 = @(10, 10)
    38 = @(5, 5)
    39 = @(9, 9)
    40 = @(1, 3)
    41 = @(1, 4)
    42 = @(7, 7)
    43 = @(6, 6)
    44 = @(5, 5)
    45 = @(5, 5)
    46 = @(2, 2)
    47 = @(8, 8)
    48 = @(4, 4)
    49 = @(1, 1)
    50 = @(6, 6)
    51 = @(4, 5)
    52 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 9).Value = $vals[0]
    $ws.Cells.Item([int]$row, 10).Value = $vals[1]
}